# Add new word/category rows to the "words list" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows (A = word, B = category), rows 51-76.
$rows = @(
    @("BUBBLE MASK",      "Bubble Mask"),
    @("Iense Gel Kajal",  "EK14 Kajal"),
    @("antioigeing",      "anti-ageing serum"),
    @("Iqvender",         "lavender clay mask"),
    @("nnlt",             "mint clay mask"),
    @("peel off",         "Charcoal peel off mask"),
    @("pearigiow",        "pearlglow serum"),
    @("sleleping",        "Sleeping Mask"),
    @("sleeping",         "Sleeping Mask"),
    @("sleeping mask",    "Sleeping Mask"),
    @("sleepingmask",     "Sleeping Mask"),
    @("iAsMINE",          "Hand Cream"),
    @("Ronantic",         "Hand Cream"),
    @("FLoRAL",           "Hand Cream"),
    @("Romantic",         "Hand Cream"),
    @("FLORAL",           "Hand Cream"),
    @("ROMANTIC",         "Hand Cream"),
    @("Floral",           "Hand Cream"),
    @("AMBER",            "Eyelashes"),
    @("3D STUDIOEFFECT",  "Eyelashes"),
    @("eyelashes",        "Eyelashes"),
    @("EYSBRAW",          "Eyebrow Styling Gel"),
    @("lotion",           "Gel Lotion"),
    @("nightcream",       "Night Cream"),
    @("daxcream",         "Day Cream"),
    @("daycream",         "Day Cream")
)

$startRow = 51
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Highlight the first new row (A51:B51) with a solid yellow fill, like the
# author did to flag where the newly appended data begins.
$ws.Range("A51:B51").Interior.Color = 65535

# Flag duplicate words in the new block (A51:A76) with Excel's standard
# "Duplicate Values" conditional formatting (light red fill / dark red text).
$dupRange = $ws.Range("A51:A76")
$cf = $dupRange.FormatConditions.AddUniqueValues()
$cf.DupeUnique = 1
$cf.Font.Color = 393372
$cf.Interior.Color = 13551615

# Reflect the author's final selection/scroll position in the saved view.
$null = $ws.Range("A67").Select()

Write-Output "done"
